$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text tweaks on existing rows ---

# Row 24: "création gestion de profil" -> "création de profil"
$ws.Range("C24").Value = "Création UI de création de profil, création de profil"

# Rows 27 & 28: "Logique gestion de profil" -> "Logique gestion de profil."
#               "En cours" -> "En cours."
$ws.Range("C27").Value = "Logique gestion de profil."
$ws.Range("E27").Value = "En cours."
$ws.Range("C28").Value = "Logique gestion de profil."
$ws.Range("E28").Value = "En cours."

# --- Fill in previously-empty rows 29, 30, 31 ---
# Copy formatting from similar, already-filled rows before writing values so
# the new rows pick up the same cell styles (date format on B, bordered /
# merged comment style on E:G for rows that carry a comment).

# Row 29 (comment present -> mirror row 28's C/D layout + row 24's E:G layout)
$ws.Range("B28").Copy()
$ws.Range("B29").PasteSpecial(-4122)
$ws.Range("C28:D28").Copy()
$ws.Range("C29:D29").PasteSpecial(-4122)
$ws.Range("E24:G24").Copy()
$ws.Range("E29:G29").PasteSpecial(-4122)

$ws.Range("B29").Value = 43160
$ws.Range("C29").Value = "Logique gestion de profil."
$ws.Range("D29").Value = "4h"
$ws.Range("E29").Value = "Gestion de profils fonctionnelle avec messages d'erreurs en cas de mauvaise manipulation."
$ws.Rows.Item(29).RowHeight = 31.5

# Row 30 (no comment -> mirror row 28's layout)
$ws.Range("B28").Copy()
$ws.Range("B30").PasteSpecial(-4122)
$ws.Range("C28:D28").Copy()
$ws.Range("C30:D30").PasteSpecial(-4122)

$ws.Range("B30").Value = 43160
$ws.Range("C30").Value = "Création UI Sélection des joueurs"
$ws.Range("D30").Value = "15 min"

# Row 31 (comment present -> mirror row 28's C/D layout + row 24's E:G layout)
$ws.Range("B28").Copy()
$ws.Range("B31").PasteSpecial(-4122)
$ws.Range("C28:D28").Copy()
$ws.Range("C31:D31").PasteSpecial(-4122)
$ws.Range("E24:G24").Copy()
$ws.Range("E31:G31").PasteSpecial(-4122)

$ws.Range("B31").Value = 43160
$ws.Range("C31").Value = "Logique sélection des joueurs "
$ws.Range("D31").Value = "1h"
$ws.Range("E31").Value = "En cours, les drop down list se remplissent. Il faut encore ajouter une erreur si deux même profils sont choisis"
$ws.Rows.Item(31).RowHeight = 47.25

# --- Update selection to match the saved view state ---
$ws.Range("I31").Select()
